# Se añade escenario Descarga de Poliza
# - Remove the extra data rows (3-5), keeping only the header row and a single
#   data row.
# - Update the remaining data row's "montoReembolso" value from 723 to 209,
#   keeping it stored as text (shared string) like the rest of the row.
# - Update the selection to reflect the new "next empty row" position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop rows 3 through 5 (the extra scenario rows), leaving header + row 2.
$ws.Rows("3:5").Delete() | Out-Null

# Set E2 ("montoReembolso") to "209" while preserving it as text (not a
# number) and without leaving behind any extra cell style.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "209"
$ws.Range("E2").Style = "Normal"

# Move the selection to where the next data row would be entered.
$ws.Range("A3:E7").Select() | Out-Null
